$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Jamal Murray'
$ws.Range("B2").Value = 'PG,SG'
$ws.Range("C2").Value = 'Denver Nuggets'
$ws.Range("A3").Value = 'Michael Porter Jr.'
$ws.Range("B3").Value = 'SF,PF'
$ws.Range("C3").Value = 'Denver Nuggets'
$ws.Range("A4").Value = 'Obi Toppin'
$ws.Range("B4").Value = 'PF'
$ws.Range("C4").Value = 'Indiana Pacers'
$ws.Range("A5").Value = 'Myles Turner'
$ws.Range("B5").Value = 'C'
$ws.Range("C5").Value = 'Indiana Pacers'
$ws.Range("A6").Value = 'Domantas Sabonis'
$ws.Range("B6").Value = 'C'
$ws.Range("C6").Value = 'Sacramento Kings'
$ws.Range("A7").Value = 'Deandre Ayton'
$ws.Range("B7").Value = 'C'
$ws.Range("C7").Value = 'Portland Trail Blazers'
$ws.Range("A8").Value = 'De''Andre Hunter'
$ws.Range("B8").Value = 'SF,PF'
$ws.Range("C8").Value = 'Atlanta Hawks'
$ws.Range("A9").Value = 'Josh Hart'
$ws.Range("B9").Value = 'SG,SF,PF'
$ws.Range("C9").Value = 'New York Knicks'
$ws.Range("A10").Value = 'Bradley Beal'
$ws.Range("B10").Value = 'PG,SG,SF'
$ws.Range("C10").Value = 'Phoenix Suns'
$ws.Range("A11").Value = 'Payton Pritchard'
$ws.Range("B11").Value = 'PG'
$ws.Range("C11").Value = 'Boston Celtics'
$ws.Range("A12").Value = 'Malik Beasley'
$ws.Range("B12").Value = 'SG,SF'
$ws.Range("C12").Value = 'Detroit Pistons'
$ws.Range("A13").Value = 'Dyson Daniels'
$ws.Range("B13").Value = 'PG,SG,SF'
$ws.Range("C13").Value = 'Atlanta Hawks'
$ws.Range("A14").Value = 'Donovan Mitchell'
$ws.Range("B14").Value = 'PG,SG'
$ws.Range("C14").Value = 'Cleveland Cavaliers'
$ws.Range("A15").Value = 'Victor Wembanyama'
$ws.Range("B15").Value = 'C'
$ws.Range("C15").Value = 'San Antonio Spurs'
$ws.Range("A16").Value = 'Kristaps Porzingis'
$ws.Range("B16").Value = 'PF,C'
$ws.Range("C16").Value = 'Boston Celtics'
$ws.Range("A17").Value = 'Cam Thomas'
$ws.Range("B17").Value = 'SG,SF'
$ws.Range("C17").Value = 'Brooklyn Nets'
$ws.Range("A18").Value = 'Tari Eason'
$ws.Range("B18").Value = 'SF,PF'
$ws.Range("C18").Value = 'Houston Rockets'
